# Edit ItemExcel.xlsx to match target revision:
#  - ItemDatas (sheet1): fix isCrafting flag (J18:J20 TRUE -> FALSE)
#  - Stats (sheet2): add new "hunger" column (H) and several new stat rows (8-20)
#  - Recipe (sheet3): no data changes, just selection/active-sheet bookkeeping
#  - Make ItemDatas the active/selected sheet (was Recipe)

$wb = $excel.ActiveWorkbook

$wsItems = $wb.Worksheets.Item("ItemDatas")
$wsStats = $wb.Worksheets.Item("Stats")
$wsRecipe = $wb.Worksheets.Item("Recipe")

# ---------------------------------------------------------------------------
# ItemDatas: isCrafting (column J) flips from TRUE to FALSE for rows 18-20
# ---------------------------------------------------------------------------
$wsItems.Range("J18").Value = $false
$wsItems.Range("J19").Value = $false
$wsItems.Range("J20").Value = $false

# ---------------------------------------------------------------------------
# Stats: new "hunger" column H, plus new rows 8-20
# ---------------------------------------------------------------------------
$wsStats.Range("H1").Value = "hunger"
$wsStats.Range("H2").Value = 0
$wsStats.Range("H3").Value = 0
$wsStats.Range("H4").Value = 0
$wsStats.Range("H5").Value = 0
$wsStats.Range("H6").Value = 0
$wsStats.Range("H7").Value = 0

$statsRows = @(
    @(2005, 0, 0, 0, 0, 0, 0, 5),
    @(2006, 0, 0, 0, 0, 0, 0, 10),
    @(2007, 0, 0, 0, 0, 0, 0, 15),
    @(2008, 0, 0, 0, 0, 0, 0, 10),
    @(3002, 50, 0, 0, 0, 0, 0, 30),
    @(3003, 0, 0, 0, 0, 0, 0, 20),
    @(3004, 0, 0, 0, 0, 0, 0, 50),
    @(3006, 20, 0, 0, 0, 0, 0, 20),
    @(3007, 0, 0, 0, 0, 0, 0, 100),
    @(3008, 0, 0, 0, 0, 0, 0, 50),
    @(3009, 0, 0, 0, 0, 0, 0, 70),
    @(4101, 0, 0, 30, 0, 0, 0, 0),
    @(6001, 100, 0, 0, 0, 0, 0, 0)
)

$r = 8
foreach ($row in $statsRows) {
    $wsStats.Range("A$r").Value = $row[0]
    $wsStats.Range("B$r").Value = $row[1]
    $wsStats.Range("C$r").Value = $row[2]
    $wsStats.Range("D$r").Value = $row[3]
    $wsStats.Range("E$r").Value = $row[4]
    $wsStats.Range("F$r").Value = $row[5]
    $wsStats.Range("G$r").Value = $row[6]
    $wsStats.Range("H$r").Value = $row[7]
    $r++
}

# narrow column H (was sized for the old bestFit text, no longer needed)
$wsStats.Columns.Item(8).ColumnWidth = 10.87

# ---------------------------------------------------------------------------
# Selections / active sheet bookkeeping
# ---------------------------------------------------------------------------
$wsStats.Activate()
$wsStats.Range("J18").Select()

$wsRecipe.Activate()
$wsRecipe.Range("F1").Select()

$wsItems.Activate()
$wsItems.Range("B8").Select()
